# Atualização de bases das ligas, do dia: 2024-02-11 às 04:25
#
# The source data for several fixtures got re-keyed; this re-syncs the
# affected rows' stat columns (B:AC) against the corrected source order
# and drops the placeholder fixture row that hadn't been played yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-content swaps: destination row => source row it should take its
# B:AC (everything except the running index in column A) content from.
# Captured as independent permutation cycles so we can snapshot every
# source row's values *before* writing any of them.
$cycles = @(
    @(51, 54),
    @(64, 65),
    @(136, 139, 137, 143, 144),
    @(138, 142),
    @(209, 210),
    @(281, 282)
)

foreach ($cycle in $cycles) {
    # Snapshot the current B:AC values for every row in this cycle first.
    $snapshot = @{}
    foreach ($r in $cycle) {
        $snapshot[$r] = $ws.Range("B$r`:AC$r").Value2()
    }

    # dest gets the snapshot belonging to the *next* row in the cycle list
    # (matches how the rows were found to have traded places).
    $n = $cycle.Count
    for ($i = 0; $i -lt $n; $i++) {
        $dest = $cycle[$i]
        $src = $cycle[($i + 1) % $n]
        $ws.Range("B$dest`:AC$dest").Value2 = $snapshot[$src]
    }
}

# The future/unplayed fixture row (id 313, match id 7785401) was removed
# from the export entirely.
$ws.Range("A315").EntireRow.Delete()
